$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix contributor name casing (re-typed values -> shared strings get re-written/reordered)
$ws.Range("C2").Value = "Hemang Jayesh S"
$ws.Range("C3").Value = "Himanshu  Bharg"
$ws.Range("C4").Value = "Saureendesai"
$ws.Range("C5").Value = "Prajapati Dipti"
$ws.Range("C8").Value = "Ranjeet Sharma"

# 2. Apply a rupee currency number format to the Spend/Contribution/Balance columns
$ws.Range("D2:F9").NumberFormat = '"₹"\ #,##0'

# 3. Drop the empty placeholder cell in G8
$ws.Range("G8").Clear()

# 4. Wipe out the extra (unused) rows 10-12 entirely
$ws.Range("A10:F12").Clear()

# Re-apply minimal formatting so the row placeholders stay but carry the right look
# -4131 = xlLeft, -4160 = xlTop
$ws.Range("A10:A12").HorizontalAlignment = -4131
$ws.Range("A10:A12").VerticalAlignment = -4160

$ws.Range("B10:B12").NumberFormat = "m/d/yyyy"

$ws.Range("F10:F12").NumberFormat = "#,##0"

# 5. Resize the columns; add a new (currently unused) column G sized like the Contributors column
$ws.Columns.Item(4).ColumnWidth = 7.1667
$ws.Columns.Item(5).ColumnWidth = 12.1667
$ws.Columns.Item(6).ColumnWidth = 8.7369
$ws.Columns.Item(7).ColumnWidth = 17.5925

# 6. Move the selection to the new column G
[void]$ws.Range("G2:G9").Select()

Write-Host "edit complete"
